# Scheduled data refresh: update market-price derived columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets with the
# latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 33141.75
$ws.Range("J3").Value = 33141.75
$ws.Range("L3").Value = 33141.75
$ws.Range("N3").Value = -33369.75

$ws.Range("H64").Value = 3666.6316
$ws.Range("I64").Value = 3427.535
$ws.Range("J64").Value = 3978.182
$ws.Range("K64").Value = 3427.535
$ws.Range("L64").Value = 3978.182
$ws.Range("M64").Value = -3179.535
$ws.Range("N64").Value = -4474.182

$ws.Range("H67").Value = 3666.6316
$ws.Range("I67").Value = 3427.535
$ws.Range("J67").Value = 3978.182
$ws.Range("K67").Value = 3427.535
$ws.Range("L67").Value = 3978.182
$ws.Range("M67").Value = -2569.535
$ws.Range("N67").Value = -5694.182

$ws.Range("H76").Value = 3420.5715
$ws.Range("I76").Value = 2986.5
$ws.Range("K76").Value = 2986.5
$ws.Range("M76").Value = -2671.5

$ws.Range("H79").Value = 3420.5715
$ws.Range("I79").Value = 2986.5
$ws.Range("K79").Value = 2986.5
$ws.Range("M79").Value = -1894.5

$ws.Range("H102").Value = 33141.75
$ws.Range("J102").Value = 33141.75
$ws.Range("L102").Value = 33141.75
$ws.Range("N102").Value = -39631.75

$ws.Range("H106").Value = 1724.091
$ws.Range("I106").Value = 1506.5
$ws.Range("K106").Value = 1506.5
$ws.Range("M106").Value = -875.5

$ws.Range("H141").Value = 1950
$ws.Range("I141").Value = 1900
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 5700
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = -520
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11638383
$ws.Range("I32").Value = 14290859
$ws.Range("J32").Value = 33801.688
$ws.Range("K32").Value = 14290859
$ws.Range("L32").Value = 33801.688
$ws.Range("M32").Value = -14290572
$ws.Range("N32").Value = -34375.688

$ws.Range("H63").Value = 1887.8889
$ws.Range("I63").Value = 1573
$ws.Range("J63").Value = 2990
$ws.Range("K63").Value = 1573
$ws.Range("L63").Value = 2990
$ws.Range("M63").Value = -887
$ws.Range("N63").Value = -4362

$ws.Range("H66").Value = 1887.8889
$ws.Range("I66").Value = 1573
$ws.Range("J66").Value = 2990
$ws.Range("K66").Value = 7865
$ws.Range("L66").Value = 14950
$ws.Range("M66").Value = -4433
$ws.Range("N66").Value = -21814

$ws.Range("H88").Value = 2047.7778
$ws.Range("I88").Value = 1982.5
$ws.Range("J88").Value = 2100
$ws.Range("K88").Value = 1982.5
$ws.Range("L88").Value = 2100
$ws.Range("M88").Value = -1576.5
$ws.Range("N88").Value = -2912

$ws.Range("H91").Value = 2047.7778
$ws.Range("I91").Value = 1982.5
$ws.Range("J91").Value = 2100
$ws.Range("K91").Value = 1982.5
$ws.Range("L91").Value = 2100
$ws.Range("M91").Value = -578.5
$ws.Range("N91").Value = -4908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1717.5312
$ws.Range("I105").Value = 1188.5
$ws.Range("J105").Value = 2129
$ws.Range("K105").Value = 1188.5
$ws.Range("L105").Value = 2129
$ws.Range("M105").Value = 558.5
$ws.Range("N105").Value = -5623

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1837.5
$ws.Range("I31").Value = 1354.1333
$ws.Range("J31").Value = 4944.857
$ws.Range("K31").Value = 1354.1333
$ws.Range("L31").Value = 4944.857
$ws.Range("M31").Value = -1059.1333
$ws.Range("N31").Value = -5534.857

$ws.Range("H34").Value = 1837.5
$ws.Range("I34").Value = 1354.1333
$ws.Range("J34").Value = 4944.857
$ws.Range("K34").Value = 1354.1333
$ws.Range("L34").Value = 4944.857
$ws.Range("M34").Value = -1152.1333
$ws.Range("N34").Value = -5348.857

$ws.Range("H62").Value = 2413
$ws.Range("I62").Value = 2400.8
$ws.Range("J62").Value = 2433.3333
$ws.Range("K62").Value = 2400.8
$ws.Range("L62").Value = 2433.3333
$ws.Range("M62").Value = -1776.8
$ws.Range("N62").Value = -3681.3333

$ws.Range("H65").Value = 2413
$ws.Range("I65").Value = 2400.8
$ws.Range("J65").Value = 2433.3333
$ws.Range("K65").Value = 12004
$ws.Range("L65").Value = 12166.6665
$ws.Range("M65").Value = -8884
$ws.Range("N65").Value = -18406.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6382.9565
$ws.Range("I70").Value = 7064.647
$ws.Range("J70").Value = 4451.5
$ws.Range("K70").Value = 7064.647
$ws.Range("L70").Value = 4451.5
$ws.Range("M70").Value = -6794.647
$ws.Range("N70").Value = -4991.5

$ws.Range("H73").Value = 6382.9565
$ws.Range("I73").Value = 7064.647
$ws.Range("J73").Value = 4451.5
$ws.Range("K73").Value = 7064.647
$ws.Range("L73").Value = 4451.5
$ws.Range("M73").Value = -6128.647
$ws.Range("N73").Value = -6323.5

$ws.Range("H80").Value = 2770.2046
$ws.Range("I80").Value = 2685.318
$ws.Range("J80").Value = 2855.0908
$ws.Range("K80").Value = 2685.318
$ws.Range("L80").Value = 2855.0908
$ws.Range("M80").Value = -1687.318
$ws.Range("N80").Value = -4851.0908

$ws.Range("H83").Value = 2770.2046
$ws.Range("I83").Value = 2685.318
$ws.Range("J83").Value = 2855.0908
$ws.Range("K83").Value = 13426.59
$ws.Range("L83").Value = 14275.454
$ws.Range("M83").Value = -8434.59
$ws.Range("N83").Value = -24259.454

$ws.Range("H105").Value = 69950
$ws.Range("J105").Value = 69950
$ws.Range("L105").Value = 69950
$ws.Range("N105").Value = -76938

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2847
$ws.Range("I82").Value = 2334.4707
$ws.Range("J82").Value = 3718.3
$ws.Range("K82").Value = 2334.4707
$ws.Range("L82").Value = 3718.3
$ws.Range("M82").Value = -1973.4707
$ws.Range("N82").Value = -4440.3

$ws.Range("H85").Value = 2847
$ws.Range("I85").Value = 2334.4707
$ws.Range("J85").Value = 3718.3
$ws.Range("K85").Value = 2334.4707
$ws.Range("L85").Value = 3718.3
$ws.Range("M85").Value = -1086.4707
$ws.Range("N85").Value = -6214.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 5000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5000
$ws.Range("N17").Value = -5344
$ws.Range("M17").ClearContents()
